$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "2025/12/03 22:00"
$ws.Range("B30").Value = "-"
$ws.Range("C30").Value = "-"
$ws.Range("D30").Value = "-"
$ws.Range("E30").Value = "-"
$ws.Range("F30").Value = "-"
$ws.Range("G30").Value = "-"
